$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet/tab
$ws.Name = "Through 2021-10-01"

# 2. Update the September row label (row 11) to drop the "(through 09-30)" suffix
$ws.Cells.Item(11, 1).Value = "September"

# 3. Insert a new row at position 12 for October data (pushes current "Total" row to 13)
$ws.Rows.Item(12).Insert()

# The Insert() call auto-populates every cell that sits under a styled column
# (D, G, M, P, S, V use the percent column style) with an empty, style-only
# cell. The source data only keeps such a placeholder for J12 (which gets a
# real value below), so clear the rest so they disappear entirely.
$ws.Cells.Item(12, 4).Clear()
$ws.Cells.Item(12, 7).Clear()
$ws.Cells.Item(12, 13).Clear()
$ws.Cells.Item(12, 16).Clear()
$ws.Cells.Item(12, 19).Clear()
$ws.Cells.Item(12, 22).Clear()

# Copy the formatting (bold + border style) from the September row label cell
# onto the new October row label cell.
$ws.Cells.Item(11, 1).Copy()
$ws.Cells.Item(12, 1).PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# 4. Populate the new October row (row 12)
$ws.Cells.Item(12, 1).Value = "October (through 10-01)"
$ws.Cells.Item(12, 8).Value = 1
$ws.Cells.Item(12, 9).Value = 5
$ws.Cells.Item(12, 10).Value = 0.1667
$ws.Cells.Item(12, 12).Value = 3
$ws.Cells.Item(12, 15).Value = 1
$ws.Cells.Item(12, 18).Value = 4
$ws.Cells.Item(12, 21).Value = 8

# 5. Update the Total row (now row 13) with the new cumulative totals
$ws.Cells.Item(13, 8).Value = 51
$ws.Cells.Item(13, 9).Value = 582
$ws.Cells.Item(13, 10).Value = 0.0806
$ws.Cells.Item(13, 12).Value = 490
$ws.Cells.Item(13, 13).Value = 0.1107
$ws.Cells.Item(13, 15).Value = 380
$ws.Cells.Item(13, 16).Value = 0.1017
$ws.Cells.Item(13, 18).Value = 852
$ws.Cells.Item(13, 19).Value = 0.0586
$ws.Cells.Item(13, 21).Value = 1178
$ws.Cells.Item(13, 22).Value = 0.0621

# 6. Narrow column A slightly (COM ColumnWidth is quantized to pixel units, so
# the closest reachable character-width value is used to approximate 23.7109375)
$ws.Columns.Item(1).ColumnWidth = 22.76
